## Day 1 Questions completed
## Adds the next set of Day-1 SDE-sheet questions (rows 6-9), marks the
## completed questions (rows 4-9) with a green "done" cell in column C,
## and widens column B to fit the longer links.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen column B for the longer links that are about to be added ---
# (86.5703125 is the exact OOXML char-width target; this host's pixel
#  grid for ColumnWidth snaps to the nearest 1/6-char step, so 85.6 is
#  the input that lands closest to it: 86.5 chars.)
$ws.Columns("B").ColumnWidth = 85.6

# --- New question rows -------------------------------------------------
# Row 6: Repeat and Missing Number
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("B6").Value = "https://www.geeksforgeeks.org/find-a-repeating-and-a-missing-number/"
$ws.Range("A6").Value = "Repeat and Missing Number "

# Row 7: Merge two sorted Arrays without extra space
$ws.Range("A5").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = "Merge two sorted Arrays without extra space "
$ws.Range("B7").Value = "https://www.geeksforgeeks.org/efficiently-merging-two-sorted-arrays-with-o1-extra-space/"

# Row 8: Kadane's algorithm
$ws.Range("A4").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("B8").Value = "https://www.geeksforgeeks.org/largest-sum-contiguous-subarray/"
$ws.Range("A8").Value = "Kadane's algorithm"

# Row 9: Merge Overlapping Subintervals
$ws.Range("A5").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value = "Merge Overlapping Subintervals"
$ws.Range("B9").Value = "https://leetcode.com/problems/merge-intervals/"

# --- Mark questions done with a green fill in column C -----------------
$ws.Range("C4:C9").Interior.Color = 5287936

# --- Selection matches the last-touched cell ----------------------------
$ws.Range("C9").Select()
